$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18 and 19 are blank template rows (only column A has a sequence
# number). Clone the formatting from the previous filled-in rows (16:17)
# so the new rows pick up the existing named styles (s=6/7/8) instead of
# the "empty slot" styles (s=10/11/12), then fill in the two new sales.
$ws.Range("B16:G17").Copy()
$ws.Range("B18:G19").PasteSpecial(-4122)

# Row 18 - VENDA 14 (17/10)
$ws.Range("B18").Value = "SEBASTIÃO RIBEIRO"
$ws.Range("C18").Value = "3c37cf19c41b99b64802f52125eef859"
$ws.Range("D18").Value = 44851
$ws.Range("E18").Value = 365
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "VENDA 14 (17/10)"

# Row 19 - VENDA 15 (17/10)
$ws.Range("B19").Value = "ALINE BASIOTE"
$ws.Range("C19").Value = "6af5116ffbe72b8492b9aa4a8b6e589e"
$ws.Range("D19").Value = 44851
$ws.Range("E19").Value = 365
$ws.Range("F19").Value = "-"
$ws.Range("G19").Value = "VENDA 15 (17/10)"
